{"js": "// Remove the student ID (\"\u5b66\u53f7\uff1a\" value) from the title block.\n// Before: \u59d3\u540d\uff1a\u5218\u6b23\u8c6a   \u5b66\u53f7\uff1a2020112921   \u5b66\u9662\u4e13\u4e1a\uff1a\u4ea4\u901a\u8fd0\u8f93\u7c7b\n// After : \u59d3\u540d\uff1a\u5218\u6b23\u8c6a   \u5b66\u53f7\uff1a            \u5b66\u9662\u4e13\u4e1a\uff1a\u4ea4\u901a\u8fd0\u8f93\u7c7b\nconst results = context.document.body.search(\"2020112921\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove the student ID (\"\u5b66\u53f7\uff1a\" value) from the title block.\n# Before: \u59d3\u540d\uff1a\u5218\u6b23\u8c6a   \u5b66\u53f7\uff1a2020112921   \u5b66\u9662\u4e13\u4e1a\uff1a\u4ea4\u901a\u8fd0\u8f93\u7c7b\n# After : \u59d3\u540d\uff1a\u5218\u6b23\u8c6a   \u5b66\u53f7\uff1a            \u5b66\u9662\u4e13\u4e1a\uff1a\u4ea4\u901a\u8fd0\u8f93\u7c7b\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2020112921\"\n$find.Replacement.Text = \"\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
